# Update "Línea 141" horarios workbook with the latest scrape (04:46:05).
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "LP1912" (main schedule sheet)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 04:46:05"
$ws1.Range("A3").Value = "Total filas: 15"

$ws1.Range("A6").Value = "04:46:05"
$ws1.Range("D6").Value = 1

$ws1.Range("A7").Value = "04:46:05"
$ws1.Range("D7").Value = 7

$ws1.Range("A8").Value = "04:46:05"
$ws1.Range("B8").Value = "05:12"
$ws1.Range("D8").Value = 26

$ws1.Range("A9").Value = "04:46:05"
$ws1.Range("D9").Value = 36

$ws1.Range("A10").Value = "04:46:05"
$ws1.Range("D10").Value = 46

$ws1.Range("A11").Value = "04:46:05"
$ws1.Range("D11").Value = 58

$ws1.Range("A12").Value = "04:46:05"
$ws1.Range("D12").Value = 66

$ws1.Range("A13").Value = "04:46:05"
$ws1.Range("D13").Value = 75

$ws1.Range("A14").Value = "04:46:05"
$ws1.Range("D14").Value = 78

$ws1.Range("A15").Value = "04:46:05"
$ws1.Range("D15").Value = 85

$ws1.Range("A16").Value = "04:46:05"
$ws1.Range("D16").Value = 98

$ws1.Range("A17").Value = "04:46:05"
$ws1.Range("D17").Value = 101

$ws1.Range("A18").Value = "04:46:05"
$ws1.Range("D18").Value = 105

$ws1.Range("A19").Value = "04:46:05"
$ws1.Range("D19").Value = 105

# New row appended to the schedule
$ws1.Range("A20").Value = "04:46:05"
$ws1.Range("B20").Value = "06:39"
$ws1.Range("C20").Value = "225_C ROCA-H SUR"
$ws1.Range("D20").Value = 113
$ws1.Range("E20").Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 04:46:05"

$ws2.Range("A6").Value = "04:46:05"
$ws2.Range("D6").Value = 1

$ws2.Range("A7").Value = "04:46:05"
$ws2.Range("D7").Value = 85

# ---------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 04:46:05"
